# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column before the
# existing "Late" column (column N), pushing "Late" / the spacer "heading"
# column / "Outstanding" one column to the right (N->O, O->P, P->Q). The new
# column inherits the column width of its left neighbour ("In Advance",
# column M). Then make "Repayment schedule" the active sheet/tab (instead of
# "Transactions") and leave the new selection on cell S7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N - existing N ("Late"), O (spacer "heading") and
# P ("Outstanding") data all shift one column to the right.
$inAdvanceColumnWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $inAdvanceColumnWidth

# Switch the active tab from "Transactions" to "Repayment schedule" and park
# the selection on S7.
$ws.Activate() | Out-Null
$ws.Range("S7").Select() | Out-Null
